# Daily attendance processing - normalize the "Recorded By" (column G) value
# so that the automated "System" entry is listed last among the recorders,
# while a literal lowercase "system" marker (if present) stays pinned first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Transform-Recorder($val) {
    $parts = $val -split ", "
    $n = $parts.Length
    if ($n -lt 2) {
        return $val
    }

    $startIdx = 0
    $prefix = ""
    $firstChar = $parts[0].Substring(0, 1)
    $firstCode = [int][char]$firstChar
    if ($firstCode -eq 115 -and $parts[0] -eq "system") {
        # Keep a leading literal "system" token pinned in place.
        $startIdx = 1
        $prefix = $parts[0] + ", "
    }

    $result = ""
    for ($i = $n - 1; $i -ge $startIdx; $i--) {
        if ($result -ne "") {
            $result = $result + ", "
        }
        $result = $result + $parts[$i]
    }

    return $prefix + $result
}

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($val -ne "") {
        $newVal = Transform-Recorder $val
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
